$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.084.56"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.304.94"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'300.33"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'98.13"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("E7").Value = "  +2.77%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D10").Value = "'36.11"
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("D11").Value = "'0.0792"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'17.69"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "2.661.49"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "2.294.00"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "'0.790"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "42.960.90"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'12.75"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "0.0₃0913"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'6.16"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'68.40"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'237.98"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "'25.10"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "'164.22"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("E30").Value = "  -13.13%  "
$ws.Range("D31").Value = "'9.14"
$ws.Range("D32").Value = "'33.15"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("D35").Value = "'4.81"
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'0.0699"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").Value = "2.020.73"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").Value = "'2.22"
$ws.Range("E45").Value = "  -4.30%  "
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").Value = "'17.60"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("D50").Value = "2.525.75"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -1.38%  "
